$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel;
# force them to stay plain text (matches the original t="inlineStr" cells),
# then drop the temporary "Text" number-format back to Normal so no stray
# cell style lingers on the cell itself.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "26.916.08"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "1.549.89"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  -0.12%  "
Set-TextValue $ws.Range("D5") "206.69"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  -0.13%  "
Set-TextValue $ws.Range("D8") "22.13"
$ws.Range("E8").Value = "  +2.98%  "
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "1.770.68"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "1.550.15"
$ws.Range("E13").Value = "  -0.18%  "
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "26.908.03"
$ws.Range("E16").Value = "  -0.19%  "
Set-TextValue $ws.Range("D17") "61.65"
$ws.Range("E17").Value = "  +0.05%  "
Set-TextValue $ws.Range("D18") "217.56"
$ws.Range("E18").Value = "  +1.23%  "
$ws.Range("E19").Value = "  +1.56%  "
Set-TextValue $ws.Range("D20") "7.26"
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  +0.36%  "
Set-TextValue $ws.Range("D23") "9.20"
$ws.Range("E23").Value = "  -0.06%  "
Set-TextValue $ws.Range("D24") "1.95"
$ws.Range("E24").Value = "  +0.38%  "
Set-TextValue $ws.Range("D25") "154.40"
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("E26").Value = "  -0.57%  "
Set-TextValue $ws.Range("D27") "14.92"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D33").Value = "1.417.98"
$ws.Range("E33").Value = "  +3.40%  "
$ws.Range("E34").Value = "  +4.18%  "
$ws.Range("E35").Value = "  +1.88%  "
Set-TextValue $ws.Range("D36") "0.972"
$ws.Range("E36").Value = "  +0.56%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  +0.46%  "
Set-TextValue $ws.Range("D39") "0.522"
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("E41").Value = "  +4.95%  "
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("E43").Value = "  +4.67%  "
Set-TextValue $ws.Range("D44") "0.992"
$ws.Range("E44").Value = "  +0.53%  "
Set-TextValue $ws.Range("D45") "64.32"
$ws.Range("E45").Value = "  +0.89%  "
$ws.Range("D47").Value = "1.684.11"
$ws.Range("E47").Value = "  -0.12%  "
Set-TextValue $ws.Range("D48") "87.66"
$ws.Range("E49").Value = "  +2.14%  "
$ws.Range("E50").Value = "  +3.63%  "
Set-TextValue $ws.Range("D51") "0.0954"
$ws.Range("E51").Value = "  +0.05%  "
